# Auto-generated update of leve profit/price data cells per upstream diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 14890.4
$ws.Range("I34").Value = 14890.4
$ws.Range("K34").Value = 14890.4
$ws.Range("M34").Value = -14687.4
$ws.Range("H36").Value = 14890.4
$ws.Range("I36").Value = 14890.4
$ws.Range("K36").Value = 14890.4
$ws.Range("M36").Value = -14175.4
$ws.Range("H88").Value = 2100
$ws.Range("I88").Value = 1285.7142
$ws.Range("K88").Value = 1285.7142
$ws.Range("M88").Value = -879.7141999999999
$ws.Range("H91").Value = 2100
$ws.Range("I91").Value = 1285.7142
$ws.Range("K91").Value = 1285.7142
$ws.Range("M91").Value = 118.2858000000001
$ws.Range("H125").Value = 111777
$ws.Range("I125").Value = 200572.6
$ws.Range("J125").Value = 782.5
$ws.Range("K125").Value = 1805153.4
$ws.Range("L125").Value = 7042.5
$ws.Range("M125").Value = -1802693.4
$ws.Range("N125").Value = -11962.5
$ws.Range("H131").Value = 6410.967
$ws.Range("I131").Value = 921.9
$ws.Range("J131").Value = 9155.5
$ws.Range("K131").Value = 2765.7
$ws.Range("L131").Value = 27466.5
$ws.Range("M131").Value = 2274.3
$ws.Range("N131").Value = -37546.5
$ws.Range("H141").Value = 5158.4375
$ws.Range("I141").Value = 1733.8889
$ws.Range("J141").Value = 23651
$ws.Range("K141").Value = 5201.6667
$ws.Range("L141").Value = 70953
$ws.Range("M141").Value = -21.66669999999976
$ws.Range("N141").Value = -81313

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 296.66666
$ws.Range("J7").Value = 281.42856
$ws.Range("L7").Value = 281.42856
$ws.Range("N7").Value = -507.42856
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50471
$ws.Range("H30").Value = 49999
$ws.Range("J30").Value = 49999
$ws.Range("L30").Value = 49999
$ws.Range("N30").Value = -50181
$ws.Range("H68").Value = 32036.25
$ws.Range("J68").Value = 32036.25
$ws.Range("L68").Value = 32036.25
$ws.Range("N68").Value = -33534.25
$ws.Range("H71").Value = 32036.25
$ws.Range("J71").Value = 32036.25
$ws.Range("L71").Value = 96108.75
$ws.Range("N71").Value = -103596.75
$ws.Range("H96").Value = 28207
$ws.Range("J96").Value = 28207
$ws.Range("L96").Value = 28207
$ws.Range("N96").Value = -33699
$ws.Range("H99").Value = 2325.9167
$ws.Range("I99").Value = 2379
$ws.Range("J99").Value = 2166.6667
$ws.Range("K99").Value = 2379
$ws.Range("L99").Value = 2166.6667
$ws.Range("M99").Value = -881
$ws.Range("N99").Value = -5162.6667
$ws.Range("H126").Value = 2325.9167
$ws.Range("I126").Value = 2379
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 7137
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -4667
$ws.Range("N126").Value = -11440.0001
$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959
$ws.Range("H132").Value = 645867.3
$ws.Range("I132").Value = 1230133.9
$ws.Range("J132").Value = 3174
$ws.Range("K132").Value = 3690401.7
$ws.Range("L132").Value = 9522
$ws.Range("M132").Value = -3687871.7
$ws.Range("N132").Value = -14582
$ws.Range("H134").Value = 2417.7896
$ws.Range("I134").Value = 1939
$ws.Range("J134").Value = 2766
$ws.Range("K134").Value = 5817
$ws.Range("L134").Value = 8298
$ws.Range("M134").Value = -3282
$ws.Range("N134").Value = -13368

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1079
$ws.Range("I122").Value = 655.6
$ws.Range("J122").Value = 1784.6666
$ws.Range("K122").Value = 5900.400000000001
$ws.Range("L122").Value = 16061.9994
$ws.Range("M122").Value = -3450.400000000001
$ws.Range("N122").Value = -20961.9994
$ws.Range("H126").Value = 4494.4165
$ws.Range("I126").Value = 3333
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 9999
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -5059
$ws.Range("N126").Value = -23680
$ws.Range("H137").Value = 27780768
$ws.Range("I137").Value = 1566
$ws.Range("J137").Value = 47623056
$ws.Range("K137").Value = 4698
$ws.Range("L137").Value = 142869168
$ws.Range("M137").Value = 402
$ws.Range("N137").Value = -142879368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3142.8572
$ws.Range("J80").Value = 3250
$ws.Range("L80").Value = 3250
$ws.Range("N80").Value = -5246
$ws.Range("H83").Value = 3142.8572
$ws.Range("J83").Value = 3250
$ws.Range("L83").Value = 16250
$ws.Range("N83").Value = -26234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5298.3335
$ws.Range("I7").Value = 5950
$ws.Range("J7").Value = 3995
$ws.Range("K7").Value = 5950
$ws.Range("L7").Value = 3995
$ws.Range("M7").Value = -5838
$ws.Range("N7").Value = -4219
$ws.Range("H18").Value = 20000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 20000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 20000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -20344
$ws.Range("H20").Value = 8592571
$ws.Range("I20").Value = 15028000
$ws.Range("J20").Value = 12000
$ws.Range("K20").Value = 15028000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = -15027774
$ws.Range("N20").Value = -12452
$ws.Range("H22").Value = 928.4286
$ws.Range("I22").Value = 1111.1111
$ws.Range("J22").Value = 599.6
$ws.Range("K22").Value = 1111.1111
$ws.Range("L22").Value = 599.6
$ws.Range("M22").Value = -816.1111000000001
$ws.Range("N22").Value = -1189.6
$ws.Range("H27").Value = 928.4286
$ws.Range("I27").Value = 1111.1111
$ws.Range("J27").Value = 599.6
$ws.Range("K27").Value = 1111.1111
$ws.Range("L27").Value = 599.6
$ws.Range("M27").Value = -1004.1111
$ws.Range("N27").Value = -813.6
$ws.Range("H46").Value = 1700
$ws.Range("I46").Value = 1800
$ws.Range("J46").Value = 1666.6666
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 1666.6666
$ws.Range("M46").Value = -1612
$ws.Range("N46").Value = -2042.6666
$ws.Range("H68").Value = 3600
$ws.Range("I68").Value = 5500
$ws.Range("J68").Value = 2650
$ws.Range("K68").Value = 5500
$ws.Range("L68").Value = 2650
$ws.Range("M68").Value = -4751
$ws.Range("N68").Value = -4148
$ws.Range("H71").Value = 3600
$ws.Range("I71").Value = 5500
$ws.Range("J71").Value = 2650
$ws.Range("K71").Value = 27500
$ws.Range("L71").Value = 13250
$ws.Range("M71").Value = -23756
$ws.Range("N71").Value = -20738
$ws.Range("H82").Value = 2040
$ws.Range("I82").Value = 1656
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1656
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1295
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 2040
$ws.Range("I85").Value = 1656
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1656
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -408
$ws.Range("N85").Value = -5496
$ws.Range("H126").Value = 5298.3335
$ws.Range("I126").Value = 5950
$ws.Range("J126").Value = 3995
$ws.Range("K126").Value = 17850
$ws.Range("L126").Value = 11985
$ws.Range("M126").Value = -15380
$ws.Range("N126").Value = -16925
$ws.Range("H133").Value = 89494.5
$ws.Range("J133").Value = 89494.5
$ws.Range("L133").Value = 89494.5
$ws.Range("N133").Value = -94554.5
$ws.Range("H136").Value = 36073904
$ws.Range("I136").Value = 52633284
$ws.Range("J136").Value = 1115212.2
$ws.Range("K136").Value = 157899852
$ws.Range("L136").Value = 3345636.6
$ws.Range("M136").Value = -157897302
$ws.Range("N136").Value = -3350736.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4530
$ws.Range("I62").Value = 3825
$ws.Range("K62").Value = 3825
$ws.Range("M62").Value = -3201
$ws.Range("H65").Value = 4530
$ws.Range("I65").Value = 3825
$ws.Range("K65").Value = 19125
$ws.Range("M65").Value = -16005
$ws.Range("H109").Value = 38299.25
$ws.Range("J109").Value = 38299.25
$ws.Range("L109").Value = 38299.25
$ws.Range("N109").Value = -41073.25
$ws.Range("H126").Value = 4628.905
$ws.Range("I126").Value = 7022.4614
$ws.Range("K126").Value = 21067.3842
$ws.Range("M126").Value = -18597.3842
